$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix city name typo: "Bankok" -> "Bangkok" (row 38, column A)
$ws.Range("A38").Value = "Bangkok"

# Update card counts (end-of-October cities update)
$ws.Range("B5").Value = 2
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = 0
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = 1
$ws.Range("B28").Value = 2
$ws.Range("B34").Value = 1
$ws.Range("B38").Value = 1
$ws.Range("B39").Value = 1
$ws.Range("B40").Value = 1
$ws.Range("B41").Value = 1
$ws.Range("B42").Value = 1
$ws.Range("B43").Value = 1
$ws.Range("B44").Value = 1
$ws.Range("B45").Value = 1

$ws.Range("A1").Select()
